$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Par_TagTechnologyToSubsets")
$ws.Activate()

# Re-apply the AutoFilter on column A (Technology) so it shows the CHP_*
# subset of technologies instead of the previous HLR_Heatpump_* subset.
$xlFilterValues = 7
$techList = @(
    "CHP_Biomass_Solid",
    "CHP_Biomass_Solid_CCS",
    "CHP_Coal_Hardcoal",
    "CHP_Coal_Hardcoal_CCS",
    "CHP_Coal_Lignite",
    "CHP_Coal_Lignite_CCS",
    "CHP_Gas_CCGT_Biogas",
    "CHP_Gas_CCGT_Biogas_CCS",
    "CHP_Gas_CCGT_Natural",
    "CHP_Gas_CCGT_Natural_CCS",
    "CHP_Gas_CCGT_SynGas",
    "CHP_Hydrogen_FuelCell",
    "CHP_Oil"
)
$ws.Range("A1:C461").AutoFilter(1, $techList, $xlFilterValues)

# Update the active cell selection left behind by the last save.
$ws.Range("A58").Select()
